# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The diff shows every data row (2-116) in column C moving from
# date serial 45188 (2023-09-19) to 45189 (2023-09-20), i.e. one day later.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45188) {
        $cell.Value2 = 45189
    }
}
